$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A508 style: drop the now-superfluous style index (3 -> default 0) ---
$ws.Range("A508").NumberFormat = "General"

# --- A521: "die Schwaeche..." entry retyped to "Arm (!= reich)" ---
$ws.Range("A521").Value = "Arm (≠ reich)"

# --- Append new vocabulary rows 532-566 ---
$ws.Range("A532").Value = "die Abbildung, -en"
$ws.Range("B532").Value = "Noun"
$ws.Range("C532").Value = "figure, illustration"
$ws.Range("D532").Value = "siehe Abbildung"

$ws.Range("A533").Value = "die Ausrede, -n"
$ws.Range("B533").Value = "Noun"
$ws.Range("C533").Value = "excuse"
$ws.Range("D533").Value = "billige Ausrede, Du hast immer eine Ausrede."

$ws.Range("A534").Value = "(sich erholen von Dat - er erholt sich – erholte sich – hat sich erholt)"
$ws.Range("B534").Value = "verb"
$ws.Range("C534").Value = "to recover"
$ws.Range("D534").Value = "von einer Überraschung sich erholen"

$ws.Range("A535").Value = "die Erholung, -en"
$ws.Range("B535").Value = "Noun"
$ws.Range("C535").Value = "The recovery, rest, relaxation"

$ws.Range("A536").Value = "der Fan, -s"
$ws.Range("B536").Value = "Noun"
$ws.Range("C536").Value = "the fan"

$ws.Range("A537").Value = "das Hallenbad, Hallenbäder"
$ws.Range("B537").Value = "Noun"
$ws.Range("C537").Value = "the indoor swimming pool"

$ws.Range("A538").Value = "das Inserat, -e"
$ws.Range("B538").Value = "Noun"
$ws.Range("C538").Value = "the advertisment"

$ws.Range("A539").Value = "die Leichtathletik"
$ws.Range("B539").Value = "Noun"
$ws.Range("C539").Value = "trach and field"

$ws.Range("A540").Value = "der Rekord, -e"
$ws.Range("B540").Value = "Noun"
$ws.Range("C540").Value = "the record"

$ws.Range("A541").Value = "Rennen- er rennt – rannte – ist gerannt"
$ws.Range("B541").Value = "verb"
$ws.Range("C541").Value = "to race, run fast"
$ws.Range("D541").Value = "er rennt die 100 meter in 11 Sekunden"

$ws.Range("A542").Value = "das Rennen, -"
$ws.Range("B542").Value = "Noun"
$ws.Range("C542").Value = "the race"

$ws.Range("A543").Value = "begrenzen- er begrenzt – begrenzte – hat begrenzt"
$ws.Range("B543").Value = "verb"
$ws.Range("C543").Value = "to limit"
$ws.Range("D543").Value = "die Geschwindigkeit auf 100km pro stunde begrenzen"

$ws.Range("A544").Value = "die Geschwindigkeit, -en"
$ws.Range("B544").Value = "Noun"
$ws.Range("C544").Value = "the speed"

$ws.Range("A545").Value = "zahlreich"
$ws.Range("B545").Value = "adj"
$ws.Range("C545").Value = "numerous"

$ws.Range("A546").Value = "nebenbei"
$ws.Range("B546").Value = "adv"
$ws.Range("C546").Value = "besides"

$ws.Range("A547").Value = "bestehen aus Dat – er besteht aus – bestand aus – hat bestanden aus"
$ws.Range("B547").Value = "verb"
$ws.Range("C547").Value = "to consist of"
$ws.Range("D547").Value = "Eine Fußballmannshaft besteht aus 11 Spielern"

$ws.Range("A548").Value = "das Resultat, -e (= das Ergebnis)"
$ws.Range("B548").Value = "Noun"
$ws.Range("C548").Value = "the result"

$ws.Range("A549").Value = "die Runde, -n"
$ws.Range("B549").Value = "Noun"
$ws.Range("C549").Value = "the lap"

$ws.Range("A550").Value = "der Sieger"
$ws.Range("B550").Value = "Noun"
$ws.Range("C550").Value = "the winner"

$ws.Range("A551").Value = "die Sportart, -en"
$ws.Range("B551").Value = "Noun"
$ws.Range("C551").Value = "the kind of sport"

$ws.Range("A552").Value = "die Sporttasche, -n"
$ws.Range("B552").Value = "Noun"
$ws.Range("C552").Value = "the sport bag"

$ws.Range("A553").Value = "der Tenisschläger, -"
$ws.Range("B553").Value = "Noun"
$ws.Range("C553").Value = "the tennis racket"

$ws.Range("A554").Value = "der Tormann, Tormänner"
$ws.Range("B554").Value = "Noun"
$ws.Range("C554").Value = "the goalkeeper"

$ws.Range("A555").Value = "das Trikot, -s"
$ws.Range("B555").Value = "Noun"
$ws.Range("C555").Value = "the jershey"

$ws.Range("A556").Value = "die Turnhalle, -n"
$ws.Range("B556").Value = "Noun"
$ws.Range("C556").Value = "the gym"

$ws.Range("A557").Value = "der Wettbewerb, -e"
$ws.Range("B557").Value = "Noun"
$ws.Range("C557").Value = "the competition"

$ws.Range("A558").Value = "Festlegen – er legt fest – legte fest – hat festgelegt"
$ws.Range("B558").Value = "verb"
$ws.Range("C558").Value = "to define, determine"

$ws.Range("A559").Value = "Siegen – er siegt – siegte – hat gesiegt"
$ws.Range("B559").Value = "verb"
$ws.Range("C559").Value = "to win"

$ws.Range("A560").Value = "Velieren – er verliert – verlor – hat verloren"
$ws.Range("B560").Value = "verb"
$ws.Range("C560").Value = "to lose"

$ws.Range("A561").Value = "Überzeugen – er überzeugt – überzeugte – hat überzeugt"
$ws.Range("B561").Value = "verb"
$ws.Range("C561").Value = "to convince"

$ws.Range("A562").Value = "Wetten – er wettet – er wettete – hat gewettet"
$ws.Range("B562").Value = "verb"
$ws.Range("C562").Value = "to bet"

$ws.Range("A563").Value = "Arm (≠ reich)"
$ws.Range("B563").Value = "adj"
$ws.Range("C563").Value = "poor"

$ws.Range("A564").Value = "die Sehenswürdigkeit, -en"
$ws.Range("B564").Value = "Noun"
$ws.Range("C564").Value = "sightseeing"

$ws.Range("A565").Value = "Bisher (= bis jetzt)"
$ws.Range("B565").Value = "adv"
$ws.Range("C565").Value = "yet"

$ws.Range("A566").Value = "Ausprobieren – er probiert aus – probierte aus – has ausprobiert"
$ws.Range("B566").Value = "verb"
$ws.Range("C566").Value = "to try out"

# --- Restore final selection/view state ---
$ws.Range("A535").Select()

Write-Host "edit complete"
